# TC20 disabled - 18NOV
# Remove the "TC20_Verify_ Find_a_Branch_Loggedinuser" sanity row (row 21)
# from the MasterExecutor sheet entirely - mirrors a manual "right click row
# header > Delete" in Excel, which removes the whole row and shifts every
# row below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 holds: Section/Page=ALL_PAGES, Functionality=END_TO_END,
# Testcase_number=TC20_Verify_ Find_a_Branch_Loggedinuser,
# Testcase_description=Verify Find a Branch for loggedinUser, RunMode=Yes,
# Severity=High. Deleting the entire row shifts rows 22:26 up to 21:25 and
# updates the sheet's used-range dimension (F26 -> F25) automatically.
$ws.Rows(21).Delete()

# Mirror the resulting selection left behind by a row-header delete: the
# (now shifted-up) row 21 remains selected, entire row.
$ws.Rows(21).Select()
